$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sach")

# Every book's TheLoai_ID (column C, rows 3-62) is bumped up by 1 - a new
# category was inserted ahead of the existing ones, shifting every
# book's category foreign key reference down by one.
for ($r = 3; $r -le 62; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value() + 1
}

# Reflect the author's final cursor position/selection on the sheet.
$ws.Activate()
$ws.Range("D56").Select()
